$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1508.7646
$ws.Range("I32").Value = 829.6
$ws.Range("J32").Value = 1791.75
$ws.Range("K32").Value = 829.6
$ws.Range("L32").Value = 1791.75
$ws.Range("M32").Value = -503.6
$ws.Range("N32").Value = -2443.75
$ws.Range("H111").Value = 15911.053
$ws.Range("I111").Value = 18959.2
$ws.Range("K111").Value = 56877.60000000001
$ws.Range("M111").Value = -53810.60000000001
$ws.Range("H132").Value = 1866.9375
$ws.Range("I132").Value = 1989.7273
$ws.Range("K132").Value = 5969.1819
$ws.Range("M132").Value = -3439.1819
$ws.Range("H138").Value = 2942.55
$ws.Range("I138").Value = 1140.8462
$ws.Range("J138").Value = 3211.77
$ws.Range("K138").Value = 3422.5386
$ws.Range("L138").Value = 9635.309999999999
$ws.Range("M138").Value = 1717.4614
$ws.Range("N138").Value = -19915.31

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 786.76746
$ws.Range("I2").Value = 728.7941
$ws.Range("J2").Value = 1005.7778
$ws.Range("K2").Value = 728.7941
$ws.Range("L2").Value = 1005.7778
$ws.Range("M2").Value = -615.7941
$ws.Range("N2").Value = -1231.7778
$ws.Range("H61").Value = 10549.771
$ws.Range("I61").Value = 8018.5
$ws.Range("J61").Value = 12681.368
$ws.Range("K61").Value = 8018.5
$ws.Range("L61").Value = 12681.368
$ws.Range("M61").Value = -7806.5
$ws.Range("N61").Value = -13105.368
$ws.Range("H63").Value = 2221.2
$ws.Range("I63").Value = 2221.2
$ws.Range("K63").Value = 2221.2
$ws.Range("M63").Value = -1535.2
$ws.Range("H66").Value = 2221.2
$ws.Range("I66").Value = 2221.2
$ws.Range("K66").Value = 11106
$ws.Range("M66").Value = -7674
$ws.Range("H74").Value = 5477.9395
$ws.Range("I74").Value = 1944
$ws.Range("J74").Value = 7244.909
$ws.Range("K74").Value = 1944
$ws.Range("L74").Value = 7244.909
$ws.Range("M74").Value = -1070
$ws.Range("N74").Value = -8992.909
$ws.Range("H77").Value = 5477.9395
$ws.Range("I77").Value = 1944
$ws.Range("J77").Value = 7244.909
$ws.Range("K77").Value = 9720
$ws.Range("L77").Value = 36224.545
$ws.Range("M77").Value = -5352
$ws.Range("N77").Value = -44960.545
$ws.Range("H97").Value = 2317.2285
$ws.Range("I97").Value = 882.35486
$ws.Range("J97").Value = 13437.5
$ws.Range("K97").Value = 882.35486
$ws.Range("L97").Value = 13437.5
$ws.Range("M97").Value = -386.35486
$ws.Range("N97").Value = -14429.5
$ws.Range("H102").Value = 100000000
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").Value = ""
$ws.Range("H111").Value = 89997.5
$ws.Range("J111").Value = 89997.5
$ws.Range("L111").Value = 89997.5
$ws.Range("N111").Value = -98177.5
$ws.Range("H112").Value = 55012.5
$ws.Range("J112").Value = 55012.5
$ws.Range("L112").Value = 55012.5
$ws.Range("N112").Value = -57966.5
$ws.Range("H113").Value = 90000
$ws.Range("J113").Value = 90000
$ws.Range("L113").Value = 90000
$ws.Range("N113").Value = -98678
$ws.Range("H114").Value = 80000
$ws.Range("J114").Value = 90000
$ws.Range("L114").Value = 90000
$ws.Range("N114").Value = -98678
$ws.Range("H116").Value = 786.76746
$ws.Range("I116").Value = 728.7941
$ws.Range("J116").Value = 1005.7778
$ws.Range("K116").Value = 728.7941
$ws.Range("L116").Value = 1005.7778
$ws.Range("M116").Value = 1565.2059
$ws.Range("N116").Value = -5593.7778
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").Value = ""
$ws.Range("H136").Value = 10549.771
$ws.Range("I136").Value = 8018.5
$ws.Range("J136").Value = 12681.368
$ws.Range("K136").Value = 24055.5
$ws.Range("L136").Value = 38044.104
$ws.Range("M136").Value = -21505.5
$ws.Range("N136").Value = -43144.104

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 786.76746
$ws.Range("I3").Value = 728.7941
$ws.Range("J3").Value = 1005.7778
$ws.Range("K3").Value = 728.7941
$ws.Range("L3").Value = 1005.7778
$ws.Range("M3").Value = -614.7941
$ws.Range("N3").Value = -1233.7778
$ws.Range("H20").Value = 1976.2046
$ws.Range("I20").Value = 1888.3438
$ws.Range("K20").Value = 1888.3438
$ws.Range("M20").Value = -1641.3438
$ws.Range("H94").Value = 3845.5908
$ws.Range("I94").Value = 1746.9445
$ws.Range("J94").Value = 13289.5
$ws.Range("K94").Value = 1746.9445
$ws.Range("L94").Value = 13289.5
$ws.Range("M94").Value = -1295.9445
$ws.Range("N94").Value = -14191.5
$ws.Range("H99").Value = 2786.2144
$ws.Range("I99").Value = 2300.889
$ws.Range("J99").Value = 3659.8
$ws.Range("K99").Value = 2300.889
$ws.Range("L99").Value = 3659.8
$ws.Range("M99").Value = -802.8890000000001
$ws.Range("N99").Value = -6655.8
$ws.Range("H105").Value = 3619.1924
$ws.Range("I105").Value = 2786.9565
$ws.Range("K105").Value = 2786.9565
$ws.Range("M105").Value = -1039.9565
$ws.Range("H107").Value = 1574.5
$ws.Range("I107").Value = 799.5
$ws.Range("K107").Value = 799.5
$ws.Range("M107").Value = 1120.5
$ws.Range("H119").Value = 87380.5
$ws.Range("J119").Value = 87380.5
$ws.Range("L119").Value = 87380.5
$ws.Range("N119").Value = -97056.5
$ws.Range("H134").Value = 6771.3
$ws.Range("I134").Value = 2955.1428
$ws.Range("K134").Value = 8865.428400000001
$ws.Range("M134").Value = -6330.428400000001

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 11832.667
$ws.Range("J55").Value = 11832.667
$ws.Range("L55").Value = 11832.667
$ws.Range("N55").Value = -12462.667

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 6345.5713
$ws.Range("I59").Value = 4999.6665
$ws.Range("K59").Value = 14998.9995
$ws.Range("M59").Value = -14458.9995
$ws.Range("H116").Value = 6055.5
$ws.Range("I116").Value = 5000
$ws.Range("J116").Value = 7111
$ws.Range("K116").Value = 15000
$ws.Range("L116").Value = 21333
$ws.Range("M116").Value = -11558
$ws.Range("N116").Value = -28217
$ws.Range("H131").Value = 313619.03
$ws.Range("J131").Value = 1764.7858
$ws.Range("L131").Value = 5294.357400000001
$ws.Range("N131").Value = -15374.3574
$ws.Range("H134").Value = 2311.2222
$ws.Range("I134").Value = 2006.8
$ws.Range("K134").Value = 6020.4
$ws.Range("M134").Value = -950.3999999999996

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 27184.666
$ws.Range("J24").Value = 27184.666
$ws.Range("L24").Value = 27184.666
$ws.Range("N24").Value = -27530.666
$ws.Range("H101").Value = 75998.5
$ws.Range("J101").Value = 75998.5
$ws.Range("L101").Value = 75998.5
$ws.Range("N101").Value = -82488.5
$ws.Range("H102").Value = 33137.41
$ws.Range("I102").Value = 3700.2
$ws.Range("J102").Value = 75190.57000000001
$ws.Range("K102").Value = 3700.2
$ws.Range("L102").Value = 75190.57000000001
$ws.Range("M102").Value = -2078.2
$ws.Range("N102").Value = -78434.57000000001
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").Value = ""
$ws.Range("H109").Value = 79000
$ws.Range("J109").Value = 79000
$ws.Range("L109").Value = 79000
$ws.Range("N109").Value = -81080
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").Value = ""
$ws.Range("H126").Value = 5399.4
$ws.Range("I126").Value = 5500
$ws.Range("J126").Value = 5374.25
$ws.Range("K126").Value = 16500
$ws.Range("L126").Value = 16122.75
$ws.Range("M126").Value = -14030
$ws.Range("N126").Value = -21062.75
$ws.Range("H132").Value = 2906.147
$ws.Range("I132").Value = 2917.1155
$ws.Range("K132").Value = 8751.3465
$ws.Range("M132").Value = -6221.3465

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 21201.77
$ws.Range("I61").Value = 21201.77
$ws.Range("K61").Value = 21201.77
$ws.Range("M61").Value = -20999.77
$ws.Range("H100").Value = 4925.5
$ws.Range("I100").Value = 5428.6924
$ws.Range("K100").Value = 5428.6924
$ws.Range("M100").Value = -4887.6924
$ws.Range("H113").Value = 21201.77
$ws.Range("I113").Value = 21201.77
$ws.Range("K113").Value = 21201.77
$ws.Range("M113").Value = -19031.77

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 450000
$ws.Range("I3").Value = 450000
$ws.Range("K3").Value = 450000
$ws.Range("M3").Value = -449886
$ws.Range("H126").Value = 1854.5428
$ws.Range("I126").Value = 1698.4667
$ws.Range("K126").Value = 5095.4001
$ws.Range("M126").Value = -2625.4001
$ws.Range("H132").Value = 1167.1637
$ws.Range("I132").Value = 1070.4147
$ws.Range("J132").Value = 1450.5
$ws.Range("K132").Value = 3211.2441
$ws.Range("L132").Value = 4351.5
$ws.Range("M132").Value = -681.2440999999999
$ws.Range("N132").Value = -9411.5
$ws.Range("H136").Value = 4196.6567
$ws.Range("I136").Value = 4558.685
$ws.Range("J136").Value = 2692.8462
$ws.Range("K136").Value = 13676.055
$ws.Range("L136").Value = 8078.5386
$ws.Range("M136").Value = -11126.055
$ws.Range("N136").Value = -13178.5386

Write-Host "Edits applied"